$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "A1"
$ws.Range("E2").Value = "DO"
$ws.Range("F2").Value = "A1"
$ws.Range("I2").Value = "M3"
$ws.Range("K2").Value = "DO"
$ws.Range("N2").Value = "M1"
$ws.Range("O2").Value = "M1"
$ws.Range("P2").Value = "M1"
$ws.Range("R2").Value = "M3"
$ws.Range("U2").Value = "A1"
$ws.Range("W2").Value = "M1"
$ws.Range("X2").Value = "A1"
$ws.Range("Y2").Value = "DO"
$ws.Range("Z2").Value = "M3"
$ws.Range("AA2").Value = "M1"
$ws.Range("AC2").Value = "M1"
$ws.Range("B3").Value = "A1"
$ws.Range("C3").Value = "M3"
$ws.Range("D3").Value = "DO"
$ws.Range("E3").Value = "M1"
$ws.Range("F3").Value = "M2"
$ws.Range("H3").Value = "M1"
$ws.Range("I3").Value = "A1"
$ws.Range("K3").Value = "A1"
$ws.Range("P3").Value = "DO"
$ws.Range("Q3").Value = "M1"
$ws.Range("R3").Value = "A2"
$ws.Range("S3").Value = "M2"
$ws.Range("T3").Value = "M3"
$ws.Range("U3").Value = "M3"
$ws.Range("V3").Value = "A1"
$ws.Range("W3").Value = "M3"
$ws.Range("X3").Value = "M2"
$ws.Range("Y3").Value = "A2"
$ws.Range("AA3").Value = "DO"
$ws.Range("AB3").Value = "M2"
$ws.Range("AC3").Value = "M1"
$ws.Range("B4").Value = "M3"
$ws.Range("E4").Value = "M1"
$ws.Range("I4").Value = "M1"
$ws.Range("K4").Value = "M1"
$ws.Range("M4").Value = "DO"
$ws.Range("N4").Value = "A1"
$ws.Range("O4").Value = "M3"
$ws.Range("P4").Value = "M1"
$ws.Range("Q4").Value = "A1"
$ws.Range("R4").Value = "M3"
$ws.Range("S4").Value = "DO"
$ws.Range("T4").Value = "M1"
$ws.Range("V4").Value = "M1"
$ws.Range("X4").Value = "M1"
$ws.Range("Z4").Value = "DO"
$ws.Range("AA4").Value = "M3"
$ws.Range("AC4").Value = "M1"
$ws.Range("B5").Value = "M2"
$ws.Range("C5").Value = "M3"
$ws.Range("D5").Value = "DO"
$ws.Range("F5").Value = "M1"
$ws.Range("G5").Value = "M2"
$ws.Range("H5").Value = "A1"
$ws.Range("L5").Value = "M3"
$ws.Range("M5").Value = "M1"
$ws.Range("N5").Value = "M1"
$ws.Range("O5").Value = "A1"
$ws.Range("P5").Value = "A1"
$ws.Range("Q5").Value = "DO"
$ws.Range("S5").Value = "M2"
$ws.Range("T5").Value = "M3"
$ws.Range("V5").Value = "M3"
$ws.Range("W5").Value = "A2"
$ws.Range("X5").Value = "M2"
$ws.Range("Y5").Value = "M1"
$ws.Range("Z5").Value = "M2"
$ws.Range("AA5").Value = "DO"
$ws.Range("AB5").Value = "M3"
$ws.Range("AC5").Value = "A1"
$ws.Range("B6").Value = "M3"
$ws.Range("D6").Value = "M1"
$ws.Range("E6").Value = "A1"
$ws.Range("F6").Value = "DO"
$ws.Range("G6").Value = "M1"
$ws.Range("I6").Value = "A2"
$ws.Range("J6").Value = "DO"
$ws.Range("M6").Value = "A1"
$ws.Range("N6").Value = "M3"
$ws.Range("P6").Value = "A2"
$ws.Range("U6").Value = "M2"
$ws.Range("W6").Value = "A2"
$ws.Range("X6").Value = "DO"
$ws.Range("Y6").Value = "M3"
$ws.Range("Z6").Value = "M1"
$ws.Range("AA6").Value = "A2"
$ws.Range("AB6").Value = "M1"
$ws.Range("AC6").Value = "A2"
$ws.Range("D7").Value = "A1"
$ws.Range("E7").Value = "M3"
$ws.Range("G7").Value = "A1"
$ws.Range("H7").Value = "A1"
$ws.Range("I7").Value = "DO"
$ws.Range("K7").Value = "A1"
$ws.Range("L7").Value = "A1"
$ws.Range("M7").Value = "M3"
$ws.Range("N7").Value = "A1"
$ws.Range("O7").Value = "A1"
$ws.Range("P7").Value = "DO"
$ws.Range("Q7").Value = "M3"
$ws.Range("R7").Value = "A1"
$ws.Range("U7").Value = "A1"
$ws.Range("Y7").Value = "A2"
$ws.Range("AB7").Value = "A2"
$ws.Range("AC7").Value = "M2"
$ws.Range("B8").Value = "DO"
$ws.Range("C8").Value = "M3"
$ws.Range("D8").Value = "A2"
$ws.Range("E8").Value = "A2"
$ws.Range("F8").Value = "A2"
$ws.Range("G8").Value = "M1"
$ws.Range("H8").Value = "M1"
$ws.Range("J8").Value = "A2"
$ws.Range("K8").Value = "A2"
$ws.Range("L8").Value = "A2"
$ws.Range("M8").Value = "DO"
$ws.Range("P8").Value = "M3"
$ws.Range("Q8").Value = "A2"
$ws.Range("R8").Value = "A2"
$ws.Range("T8").Value = "A2"
$ws.Range("U8").Value = "M1"
$ws.Range("W8").Value = "DO"
$ws.Range("X8").Value = "A1"
$ws.Range("Y8").Value = "A2"
$ws.Range("Z8").Value = "A2"
$ws.Range("AA8").Value = "A1"
$ws.Range("AC8").Value = "M2"
$ws.Range("B9").Value = "M3"
$ws.Range("C9").Value = "DO"
$ws.Range("D9").Value = "M1"
$ws.Range("E9").Value = "M2"
$ws.Range("F9").Value = "M2"
$ws.Range("G9").Value = "A1"
$ws.Range("H9").Value = "A2"
$ws.Range("I9").Value = "DO"
$ws.Range("J9").Value = "M1"
$ws.Range("K9").Value = "M3"
$ws.Range("L9").Value = "M2"
$ws.Range("M9").Value = "M1"
$ws.Range("O9").Value = "A2"
$ws.Range("P9").Value = "A1"
$ws.Range("R9").Value = "M2"
$ws.Range("S9").Value = "DO"
$ws.Range("T9").Value = "M1"
$ws.Range("V9").Value = "M3"
$ws.Range("W9").Value = "M1"
$ws.Range("Y9").Value = "DO"
$ws.Range("Z9").Value = "M1"
$ws.Range("AA9").Value = "M2"
$ws.Range("AB9").Value = "A2"
$ws.Range("AC9").Value = "A2"
$ws.Range("B10").Value = "A1"
$ws.Range("C10").Value = "M3"
$ws.Range("D10").Value = "M2"
$ws.Range("E10").Value = "DO"
$ws.Range("F10").Value = "A2"
$ws.Range("G10").Value = "M1"
$ws.Range("H10").Value = "M2"
$ws.Range("I10").Value = "M1"
$ws.Range("K10").Value = "M2"
$ws.Range("L10").Value = "DO"
$ws.Range("N10").Value = "M3"
$ws.Range("Q10").Value = "M2"
$ws.Range("R10").Value = "DO"
$ws.Range("S10").Value = "M1"
$ws.Range("U10").Value = "M3"
$ws.Range("W10").Value = "A1"
$ws.Range("X10").Value = "A2"
$ws.Range("Y10").Value = "M2"
$ws.Range("Z10").Value = "M1"
$ws.Range("AA10").Value = "M2"
$ws.Range("AC10").Value = "DO"
$ws.Range("Z3").Value = "A2"
$ws.Range("U9").Value = "M2"
$ws.Range("M10").Value = "A2"
$ws.Range("V10").Value = "A2"
